# Generate Report for Handback
# Marks the "a.md" file as handed back for both the zh-cn and de-de
# localization targets: status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the Latest Target File / Latest
# Handback File / Latest Handback DateTime columns get populated, and a
# hyperlink to a.md is added next to the new Latest Target File value.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status column updates (Overview + both language sheets) ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- zh-cn: Latest Target File / Latest Handback File / hyperlinks ---
$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8274e1c7dbb6ee85c6d364bd08b69d694c445a4a/e2e/a.md"

$zhcn.Range("I2").Value = "a.md"
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("I3").Value = "a.md"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# zh-cn Latest Handback DateTime moves from the placeholder date to a real one
$zhcn.Range("K2").Value = "2016-08-31 08:44:40"
$zhcn.Range("K3").Value = "2016-08-31 08:44:40"

# --- de-de: Latest Target File / Latest Handback File / hyperlinks ---
$dede.Range("I2").Value = "a.md"
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("I3").Value = "a.md"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

# de-de Latest Handback DateTime moves from the placeholder date to a real one
$dede.Range("K2").Value = "2016-08-31 08:44:48"
$dede.Range("K3").Value = "2016-08-31 08:44:48"

# --- Hyperlinks on the new "Latest Target File" cells ---
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $aUrl, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $aUrl, "", "", "a.md")

$dede.Hyperlinks.Add($dede.Range("I2"), $aUrl, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $aUrl, "", "", "a.md")

# --- Column width adjustments to fit the new, longer text ---
# (ColumnWidth values below are chosen so the resulting stored column
# width lands as close as possible to the target width, since the host
# quantizes ColumnWidth -> stored width in 1/6 character-unit steps.)
$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667

$zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$zhcn.Columns.Item(10).ColumnWidth = 39.1666666666667

$dede.Columns.Item(3).ColumnWidth = 29.1666666666667
$dede.Columns.Item(10).ColumnWidth = 39.1666666666667
